$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $text = $cell.Text

    if ($text -eq "") { continue }

    $parts = $text -split ", "
    $systemParts = @()
    $otherParts = @()
    foreach ($p in $parts) {
        # Case-sensitive comparison: only the exact token "System" (capital S)
        # is moved to the front; "system" (lowercase) is left in place.
        if ($p.CompareTo("System") -eq 0) {
            $systemParts += $p
        } else {
            $otherParts += $p
        }
    }

    if ($systemParts.Count -eq 0) { continue }

    $newParts = $systemParts + $otherParts
    $newText = $newParts -join ", "

    if ($newText.CompareTo($text) -ne 0) {
        $cell.Value = $newText
    }
}
